$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (Unicorn Name shifts to column C)
$ws.Columns("B").Insert()

# New header and values for the inserted "Month" column
$ws.Range("B1").Value = "Month"
$ws.Range("B2").Value = "November"
$ws.Range("B3").Value = "April"

# Clear out the old unicorn name values (now in column C), keep header
$ws.Range("C2:C3").ClearContents()

$ws.Columns("B").AutoFit()

$ws.Range("B4").Select()
